# Update the width value for the Car.png entry (row 4) in the image metadata
# table, reflecting that the image was resized/updated.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 722
